$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing header cells (G1,H1,I1) one column to the right
# (G1->H1, H1->I1, I1->J1), then drop in the new header at G1. Copying
# the values directly (rather than using Range.Insert, which would also
# shift the <cols> width definitions) matches the target, where columns
# H/I keep their original widths and only column G's width changes.
$ws.Range("J1").Value = $ws.Range("I1").Text
$ws.Range("I1").Value = $ws.Range("H1").Text
$ws.Range("H1").Value = $ws.Range("G1").Text
$ws.Range("G1").Value = "發表論文名稱"

# J1 is a brand-new cell; give it the same header formatting as the rest
# of row 1 (style index 1) without introducing a new style entry.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Adjust column width for column G to fit the new header
$ws.Columns.Item(7).ColumnWidth = 15.25

# Update the active selection as recorded in the file
$ws.Range("H7").Select()
